$p = $ppt.ActivePresentation

# Slide 2: "TextBox 3" shape - merge "The" " " "Moon" runs into a single "The Moon" run
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item("TextBox 3").TextFrame.TextRange
[void]$tr2.Delete()
[void]$s2.Shapes.Item("TextBox 3").TextFrame.TextRange.InsertAfter("The Moon")

# Slide 3: "Title 1" shape - merge "One" " " "More" runs into a single "One More" run
$s3 = $p.Slides.Item(3)
$tr3title = $s3.Shapes.Item("Title 1").TextFrame.TextRange
[void]$tr3title.Delete()
[void]$s3.Shapes.Item("Title 1").TextFrame.TextRange.InsertAfter("One More")

# Slide 3: "TextBox 3" shape - merge "The" " " "Moon" runs into a single "The Moon" run
$tr3textbox = $s3.Shapes.Item("TextBox 3").TextFrame.TextRange
[void]$tr3textbox.Delete()
[void]$s3.Shapes.Item("TextBox 3").TextFrame.TextRange.InsertAfter("The Moon")
